$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Drop the stale "_GoBack" bookmark that sat at the end of the
#    title ("Iteration Plan 3") paragraph. Word will re-number the
#    remaining bookmark ids (OLE_LINK1) automatically on save.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) In the milestones table, the "Sign in page" row's "Assigned to
#    (name)" cell currently reads "Aaron" - change it to "Charnes".
#    Locate the row by its description text ("Sign in page") instead
#    of hard-coding table/row numbers, then use column 5 ("Assigned
#    to (name)") within that row.
# ------------------------------------------------------------------
$table = $null
$targetRow = 0
for ($ti = 1; $ti -le $d.Tables.Count; $ti++) {
    $candidate = $d.Tables($ti)
    if ($candidate.Columns.Count -ge 5) {
        for ($r = 1; $r -le $candidate.Rows.Count; $r++) {
            $descCell = $candidate.Cell($r, 2)
            if ($descCell.Range.Text.TrimEnd([char]13, [char]7) -eq "Sign in page") {
                $table = $candidate
                $targetRow = $r
                break
            }
        }
    }
    if ($targetRow -ne 0) {
        break
    }
}

if ($targetRow -eq 0) {
    throw "Could not find the 'Sign in page' milestone row"
}

$cell = $table.Cell($targetRow, 5)
$cellStart = $cell.Range.Start
$nameLen = $cell.Range.Text.TrimEnd([char]13, [char]7).Length

$nameRng = $d.Range($cellStart, $cellStart + $nameLen)
if ($nameRng.Text -ne "Aaron") {
    throw "Unexpected cell content: [" + $nameRng.Text + "]"
}

# Replace "Aaron" with "Charnes" plus a throwaway marker character. The
# marker gives us a safe (non-paragraph-mark) insertion point to anchor
# the new bookmark on; collapsing a range directly onto the hidden
# end-of-cell mark is unreliable, so we add the bookmark around the
# marker first and then delete the marker through the bookmark's own
# Range, which lets the bookmark naturally collapse to sit right after
# "Charnes".
$nameRng.Text = "CharnesZ"

$markerRng = $d.Range($cellStart + 7, $cellStart + 8)
if ($markerRng.Text -ne "Z") {
    throw "Marker not where expected: [" + $markerRng.Text + "]"
}

# ------------------------------------------------------------------
# 3) Re-create "_GoBack" collapsed right after the new "Charnes" run.
# ------------------------------------------------------------------
$d.Bookmarks.Add("_GoBack", $markerRng)
$goBack = $d.Bookmarks("_GoBack")
$goBackRange = $goBack.Range
$goBackRange.Text = ""
